$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the class-count values in column B
$ws.Range("B2").Value = 2
$ws.Range("B4").Value = 4
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 6
$ws.Range("B8").Value = -1
$ws.Range("B10").Value = 0
$ws.Range("B12").Value = 7
$ws.Range("B13").Value = 3
$ws.Range("B14").Value = 1

# Move the active cell selection to E20 (mark as done)
$ws.Range("E20").Select()
